# Scheduled-runner refresh: overwrite market-snapshot columns
# (currentAveragePrice / NQ / HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ)
# on the affected leve rows across all eight crafting-job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 279.13333
$ws.Range("I33").Value = 279.13333
$ws.Range("K33").Value = 279.13333
$ws.Range("M33").Value = -50.13333
$ws.Range("H40").Value = 71430950
$ws.Range("J40").Value = 83335800
$ws.Range("L40").Value = 83335800
$ws.Range("N40").Value = -83336150
$ws.Range("H98").Value = 1848.4736
$ws.Range("I98").Value = 1520.8857
$ws.Range("K98").Value = 1520.8857
$ws.Range("M98").Value = -22.88570000000004
$ws.Range("H99").Value = 2947.7144
$ws.Range("J99").Value = 4049.6
$ws.Range("L99").Value = 12148.8
$ws.Range("N99").Value = -15144.8
$ws.Range("H122").Value = 1848.4736
$ws.Range("I122").Value = 1520.8857
$ws.Range("K122").Value = 4562.6571
$ws.Range("M122").Value = -2112.6571
$ws.Range("H138").Value = 3074.104
$ws.Range("I138").Value = 1872.7778
$ws.Range("J138").Value = 4128.927
$ws.Range("K138").Value = 5618.3334
$ws.Range("L138").Value = 12386.781
$ws.Range("M138").Value = -478.3334000000004
$ws.Range("N138").Value = -22666.781

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1733.6818
$ws.Range("I2").Value = 1647.2222
$ws.Range("J2").Value = 2122.75
$ws.Range("K2").Value = 1647.2222
$ws.Range("L2").Value = 2122.75
$ws.Range("M2").Value = -1534.2222
$ws.Range("N2").Value = -2348.75
$ws.Range("H60").Value = 78177.336
$ws.Range("J60").Value = 24999
$ws.Range("L60").Value = 24999
$ws.Range("N60").Value = -26465
$ws.Range("H61").Value = 2841069.8
$ws.Range("I61").Value = 3128023.5
$ws.Range("J61").Value = 1004566.6
$ws.Range("K61").Value = 3128023.5
$ws.Range("L61").Value = 1004566.6
$ws.Range("M61").Value = -3127811.5
$ws.Range("N61").Value = -1004990.6
$ws.Range("H74").Value = 3014.077
$ws.Range("I74").Value = 2432.0417
$ws.Range("K74").Value = 2432.0417
$ws.Range("M74").Value = -1558.0417
$ws.Range("H77").Value = 3014.077
$ws.Range("I77").Value = 2432.0417
$ws.Range("K77").Value = 12160.2085
$ws.Range("M77").Value = -7792.208500000001
$ws.Range("H116").Value = 1733.6818
$ws.Range("I116").Value = 1647.2222
$ws.Range("J116").Value = 2122.75
$ws.Range("K116").Value = 1647.2222
$ws.Range("L116").Value = 2122.75
$ws.Range("M116").Value = 646.7778000000001
$ws.Range("N116").Value = -6710.75
$ws.Range("H136").Value = 2841069.8
$ws.Range("I136").Value = 3128023.5
$ws.Range("J136").Value = 1004566.6
$ws.Range("K136").Value = 9384070.5
$ws.Range("L136").Value = 3013699.8
$ws.Range("M136").Value = -9381520.5
$ws.Range("N136").Value = -3018799.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1733.6818
$ws.Range("I3").Value = 1647.2222
$ws.Range("J3").Value = 2122.75
$ws.Range("K3").Value = 1647.2222
$ws.Range("L3").Value = 2122.75
$ws.Range("M3").Value = -1533.2222
$ws.Range("N3").Value = -2350.75
$ws.Range("H94").Value = 3553.182
$ws.Range("I94").Value = 7274.5
$ws.Range("J94").Value = 1426.7142
$ws.Range("K94").Value = 7274.5
$ws.Range("L94").Value = 1426.7142
$ws.Range("M94").Value = -6823.5
$ws.Range("N94").Value = -2328.7142

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1574.7059
$ws.Range("I16").Value = 1195
$ws.Range("K16").Value = 1195
$ws.Range("M16").Value = -908
$ws.Range("H94").Value = 1309
$ws.Range("I94").Value = 1222
$ws.Range("J94").Value = 1352.5
$ws.Range("K94").Value = 1222
$ws.Range("L94").Value = 1352.5
$ws.Range("M94").Value = -771
$ws.Range("N94").Value = -2254.5
$ws.Range("H113").Value = 1574.7059
$ws.Range("I113").Value = 1195
$ws.Range("K113").Value = 1195
$ws.Range("M113").Value = 975

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 433
$ws.Range("I7").Value = 433
$ws.Range("K7").Value = 1299
$ws.Range("M7").Value = -1187
$ws.Range("H23").Value = 1342.8182
$ws.Range("J23").Value = 1619.1111
$ws.Range("L23").Value = 4857.3333
$ws.Range("N23").Value = -5327.3333
$ws.Range("H25").Value = 22555
$ws.Range("I25").Value = 17166
$ws.Range("J25").Value = 33333
$ws.Range("K25").Value = 51498
$ws.Range("L25").Value = 99999
$ws.Range("M25").Value = -51329
$ws.Range("N25").Value = -100337
$ws.Range("H30").Value = 22555
$ws.Range("I30").Value = 17166
$ws.Range("J30").Value = 33333
$ws.Range("K30").Value = 51498
$ws.Range("L30").Value = 99999
$ws.Range("M30").Value = -51396
$ws.Range("N30").Value = -100203
$ws.Range("H97").Value = 1454
$ws.Range("J97").Value = 1205.875
$ws.Range("L97").Value = 3617.625
$ws.Range("N97").Value = -4609.625
$ws.Range("H132").Value = 1999.5
$ws.Range("I132").Value = 1819.3
$ws.Range("J132").Value = 2299.8333
$ws.Range("K132").Value = 16373.7
$ws.Range("L132").Value = 20698.4997
$ws.Range("M132").Value = -13843.7
$ws.Range("N132").Value = -25758.4997
$ws.Range("H137").Value = 6686.0835
$ws.Range("J137").Value = 13277.667
$ws.Range("L137").Value = 39833.001
$ws.Range("N137").Value = -50033.001
$ws.Range("H139").Value = 4026.5
$ws.Range("I139").Value = 1519.6428
$ws.Range("K139").Value = 4558.928400000001
$ws.Range("M139").Value = 581.0715999999993

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3238.5557
$ws.Range("I102").Value = 3030
$ws.Range("K102").Value = 3030
$ws.Range("M102").Value = -1408
$ws.Range("H126").Value = 2133.6843
$ws.Range("I126").Value = 1952.9286
$ws.Range("J126").Value = 2639.8
$ws.Range("K126").Value = 5858.7858
$ws.Range("L126").Value = 7919.400000000001
$ws.Range("M126").Value = -3388.7858
$ws.Range("N126").Value = -12859.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H45").Value = 54999
$ws.Range("I45").Value = 9999
$ws.Range("K45").Value = 9999
$ws.Range("M45").Value = -9592
$ws.Range("H115").Value = 150301.5
$ws.Range("J115").Value = 150301.5
$ws.Range("L115").Value = 150301.5
$ws.Range("N115").Value = -152651.5
$ws.Range("H132").Value = 3328.5715
$ws.Range("I132").Value = 2022.6666
$ws.Range("K132").Value = 6067.9998
$ws.Range("M132").Value = -3537.9998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 14428.571
$ws.Range("I31").Value = 14428.571
$ws.Range("K31").Value = 14428.571
$ws.Range("M31").Value = -14080.571
$ws.Range("H37").Value = 89497.5
$ws.Range("I37").Value = 99000
$ws.Range("K37").Value = 99000
$ws.Range("M37").Value = -98797
$ws.Range("H47").Value = 45000
$ws.Range("J47").Value = 45000
$ws.Range("L47").Value = 45000
$ws.Range("N47").Value = -46144
$ws.Range("H81").Value = 1972.75
$ws.Range("I81").Value = 1863.3334
$ws.Range("J81").Value = 2301
$ws.Range("K81").Value = 3726.6668
$ws.Range("L81").Value = 4602
$ws.Range("M81").Value = -2665.6668
$ws.Range("N81").Value = -6724
$ws.Range("H84").Value = 1972.75
$ws.Range("I84").Value = 1863.3334
$ws.Range("J84").Value = 2301
$ws.Range("K84").Value = 18633.334
$ws.Range("L84").Value = 23010
$ws.Range("M84").Value = -13329.334
$ws.Range("N84").Value = -33618
$ws.Range("H126").Value = 2760.6553
$ws.Range("I126").Value = 2871.4546
$ws.Range("K126").Value = 8614.363799999999
$ws.Range("M126").Value = -6144.363799999999
